$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for rows with changed values.
# NumberFormat is forced to text ("@") before assigning D-column values so that
# numeric-looking strings (e.g. "0.9974") are preserved as text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.349.14"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.67"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9974"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.97"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6263"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07492"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2902"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.37"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.07"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.995"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6788"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001032"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.09"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.101.80"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.139"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.394.98"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.00"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.33"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9986"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.461"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9976"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.56"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1374"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.403"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.51"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06465"
$ws.Range("E29").Value = "  +15.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.374"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.091"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.062"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6994"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.575"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.260.27"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.834"
$ws.Range("E39").Value = "  +4.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01826"
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.607"
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9098"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9982"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.006.38"
$ws.Range("E44").Value = "  -18.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.48"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.07"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.077"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1173"
$ws.Range("E49").Value = "  +2.70%  "

# Row 50: BabyDogeCoin dropped from list; EnergySwap moves up from row 51 with new figures.
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.049"
$ws.Range("E50").Value = "  +0.69%  "

# Row 51: TheSandbox newly added, replacing EnergySwap previously in this row.
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3941"
$ws.Range("E51").Value = "  -1.85%  "
